# Update the public EPEX Spot / Gaz / CO2 Excel workbook with the latest day of data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": add a new column AR for 27-jul, copying the header style
# of the existing last header cell (AQ1) and filling in the 24 hourly prices.
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("AQ1").Copy()
$wsPrix.Range("AR1").PasteSpecial(-4122)  # xlPasteFormats: copy the bold/centered header style
$wsPrix.Range("AR1").Value = "27-jul"

$prixValues = @(60.5, 49.86, 48.9, 32.09, 37.92, 40.01, 40.16, 50, 32.34, 23.69, 14.77, 33.17, 31.66, 15.04, 8.85, 12.93, 16.14, 24, 23.3, 35.32, 40, 62.67, 88.64, 66.95)

for ($i = 0; $i -lt $prixValues.Count; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 44).Value = $prixValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append a new row 41 for 2025-07-25.
# The date must stay stored as literal text (like every other date cell in
# column A), so force text formatting before assigning the string value,
# then drop back to the default "Normal" style (no explicit style index),
# matching the rest of column A.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A41").NumberFormat = "@"
$wsGaz.Range("A41").Value = "2025-07-25"
$wsGaz.Range("A41").Style = "Normal"
$wsGaz.Range("B41").Value = 31.7

# ---------------------------------------------------------------------------
# Sheet "CO2": append a new row 41 for 2025-07-25.
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A41").NumberFormat = "@"
$wsCo2.Range("A41").Value = "2025-07-25"
$wsCo2.Range("A41").Style = "Normal"
$wsCo2.Range("B41").Value = 70.7
